$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# New Abnormal Return (AR) values for rows 5-25 (column B), and the
# corresponding Variance AR (column C). CAR (D), Variance CAR (E),
# T-stat (F) and P-value (G) are recomputed for the aggregated event window.

$ws.Cells.Item(5, 2).Value = -0.003935650435941625
$ws.Cells.Item(5, 3).Value = 0.0001593037686820239
$ws.Cells.Item(5, 4).Value = -0.003935650435941625
$ws.Cells.Item(5, 5).Value = 0.0001593037686820239
$ws.Cells.Item(5, 6).Value = -0.3118196591099094
$ws.Cells.Item(5, 7).Value = 0.3776974515317684

$ws.Cells.Item(6, 2).Value = 0.009886099242211092
$ws.Cells.Item(6, 3).Value = 0.0001593037686820239
$ws.Cells.Item(6, 4).Value = 0.005950448806269467
$ws.Cells.Item(6, 5).Value = 0.0003186075373640478
$ws.Cells.Item(6, 6).Value = 0.3333663041389559
$ws.Cells.Item(6, 7).Value = 0.3695457225520961

$ws.Cells.Item(7, 2).Value = 0.03686949311152292
$ws.Cells.Item(7, 3).Value = 0.0001593037686820239
$ws.Cells.Item(7, 4).Value = 0.04281994191779238
$ws.Cells.Item(7, 5).Value = 0.0004779113060460718
$ws.Cells.Item(7, 6).Value = 1.958720286988042
$ws.Cells.Item(7, 7).Value = 0.02553743279552179

$ws.Cells.Item(8, 2).Value = 0.003929652349630477
$ws.Cells.Item(8, 3).Value = 0.0001593037686820239
$ws.Cells.Item(8, 4).Value = 0.04674959426742286
$ws.Cells.Item(8, 5).Value = 0.0006372150747280957
$ws.Cells.Item(8, 6).Value = 1.851973744272164
$ws.Cells.Item(8, 7).Value = 0.03250760010180109

$ws.Cells.Item(9, 2).Value = -0.01491602802465154
$ws.Cells.Item(9, 3).Value = 0.0001593037686820239
$ws.Cells.Item(9, 4).Value = 0.03183356624277132
$ws.Cells.Item(9, 5).Value = 0.0007965188434101196
$ws.Cells.Item(9, 6).Value = 1.12794329553988
$ws.Cells.Item(9, 7).Value = 0.1301240703682524

$ws.Cells.Item(10, 2).Value = -0.00794981979211278
$ws.Cells.Item(10, 3).Value = 0.0001593037686820239
$ws.Cells.Item(10, 4).Value = 0.02388374645065854
$ws.Cells.Item(10, 5).Value = 0.0009558226120921435
$ws.Cells.Item(10, 6).Value = 0.7725272397247469
$ws.Cells.Item(10, 7).Value = 0.2202062168233466

$ws.Cells.Item(11, 2).Value = -0.006848667889507258
$ws.Cells.Item(11, 3).Value = 0.0001593037686820239
$ws.Cells.Item(11, 4).Value = 0.01703507856115128
$ws.Cells.Item(11, 5).Value = 0.001115126380774168
$ws.Cells.Item(11, 6).Value = 0.5101314461737966
$ws.Cells.Item(11, 7).Value = 0.3051678688567327

$ws.Cells.Item(12, 2).Value = -0.009884378583264674
$ws.Cells.Item(12, 3).Value = 0.0001593037686820239
$ws.Cells.Item(12, 4).Value = 0.007150699977886603
$ws.Cells.Item(12, 5).Value = 0.001274430149456191
$ws.Cells.Item(12, 6).Value = 0.2003044225103632
$ws.Cells.Item(12, 7).Value = 0.4206893373338501

$ws.Cells.Item(13, 2).Value = 0.00008549523150954211
$ws.Cells.Item(13, 3).Value = 0.0001593037686820239
$ws.Cells.Item(13, 4).Value = 0.007236195209396145
$ws.Cells.Item(13, 5).Value = 0.001433733918138215
$ws.Cells.Item(13, 6).Value = 0.191106735737216
$ws.Cells.Item(13, 7).Value = 0.4242858361062187

$ws.Cells.Item(14, 2).Value = -0.004063667302045998
$ws.Cells.Item(14, 3).Value = 0.0001593037686820239
$ws.Cells.Item(14, 4).Value = 0.003172527907350148
$ws.Cells.Item(14, 5).Value = 0.001593037686820239
$ws.Cells.Item(14, 6).Value = 0.07948632644739234
$ws.Cells.Item(14, 7).Value = 0.4683494894655674

$ws.Cells.Item(15, 2).Value = 0.000899859998840269
$ws.Cells.Item(15, 3).Value = 0.0001593037686820239
$ws.Cells.Item(15, 4).Value = 0.004072387906190416
$ws.Cells.Item(15, 5).Value = 0.001752341455502263
$ws.Cells.Item(15, 6).Value = 0.09728363091826903
$ws.Cells.Item(15, 7).Value = 0.4612831664763337

$ws.Cells.Item(16, 2).Value = -0.009952276246002801
$ws.Cells.Item(16, 3).Value = 0.0001593037686820239
$ws.Cells.Item(16, 4).Value = -0.005879888339812385
$ws.Cells.Item(16, 5).Value = 0.001911645224184287
$ws.Cells.Item(16, 6).Value = -0.1344823937235355
$ws.Cells.Item(16, 7).Value = 0.4465558055837257

$ws.Cells.Item(17, 2).Value = -0.01616381859091522
$ws.Cells.Item(17, 3).Value = 0.0001593037686820239
$ws.Cells.Item(17, 4).Value = -0.02204370693072761
$ws.Cells.Item(17, 5).Value = 0.002070948992866311
$ws.Cells.Item(17, 6).Value = -0.4843953074283416
$ws.Cells.Item(17, 7).Value = 0.3142299889832101

$ws.Cells.Item(18, 2).Value = 0.008899899328894824
$ws.Cells.Item(18, 3).Value = 0.0001593037686820239
$ws.Cells.Item(18, 4).Value = -0.01314380760183279
$ws.Cells.Item(18, 5).Value = 0.002230252761548335
$ws.Cells.Item(18, 6).Value = -0.278319829957595
$ws.Cells.Item(18, 7).Value = 0.3904796133864255

$ws.Cells.Item(19, 2).Value = 0.007449022385713353
$ws.Cells.Item(19, 3).Value = 0.0001593037686820239
$ws.Cells.Item(19, 4).Value = -0.005694785216119432
$ws.Cells.Item(19, 5).Value = 0.002389556530230359
$ws.Cells.Item(19, 6).Value = -0.1164980602849768
$ws.Cells.Item(19, 7).Value = 0.4536680156309918

$ws.Cells.Item(20, 2).Value = 0.005613479187729091
$ws.Cells.Item(20, 3).Value = 0.0001593037686820239
$ws.Cells.Item(20, 4).Value = -0.00008130602839034103
$ws.Cells.Item(20, 5).Value = 0.002548860298912383
$ws.Cells.Item(20, 6).Value = -0.001610459215630958
$ws.Cells.Item(20, 7).Value = 0.4993580569717858

$ws.Cells.Item(21, 2).Value = -0.003571479312218215
$ws.Cells.Item(21, 3).Value = 0.0001593037686820239
$ws.Cells.Item(21, 4).Value = -0.003652785340608557
$ws.Cells.Item(21, 5).Value = 0.002708164067594407
$ws.Cells.Item(21, 6).Value = -0.07019184620400255
$ws.Cells.Item(21, 7).Value = 0.4720439430898959

$ws.Cells.Item(22, 2).Value = -0.008772209685333602
$ws.Cells.Item(22, 3).Value = 0.0001593037686820239
$ws.Cells.Item(22, 4).Value = -0.01242499502594216
$ws.Cells.Item(22, 5).Value = 0.00286746783627643
$ws.Cells.Item(22, 6).Value = -0.2320314991130925
$ws.Cells.Item(22, 7).Value = 0.4083361375313432

$ws.Cells.Item(23, 2).Value = 0.005283362004705934
$ws.Cells.Item(23, 3).Value = 0.0001593037686820239
$ws.Cells.Item(23, 4).Value = -0.007141633021236225
$ws.Cells.Item(23, 5).Value = 0.003026771604958455
$ws.Cells.Item(23, 6).Value = -0.1298098668088332
$ws.Cells.Item(23, 7).Value = 0.448402069746998

$ws.Cells.Item(24, 2).Value = -0.001932949177437579
$ws.Cells.Item(24, 3).Value = 0.0001593037686820239
$ws.Cells.Item(24, 4).Value = -0.009074582198673804
$ws.Cells.Item(24, 5).Value = 0.003186075373640478
$ws.Cells.Item(24, 6).Value = -0.1607676323909191
$ws.Cells.Item(24, 7).Value = 0.4361924923206937

$ws.Cells.Item(25, 2).Value = 0.0006937044841925427
$ws.Cells.Item(25, 3).Value = 0.0001593037686820239
$ws.Cells.Item(25, 4).Value = -0.00838087771448126
$ws.Cells.Item(25, 5).Value = 0.003345379142322502
$ws.Cells.Item(25, 6).Value = -0.1448994815198182
$ws.Cells.Item(25, 7).Value = 0.4424439297763026

